# Mise à jour de l'application
# Adds a new training-date column (BX) for 2025-11-07 (serial 45968),
# copying the BW column's formatting (number/cell styles) into BX, then
# filling in each player's attendance mark for the new date.
# Row 12 has no data beyond column AX (player inactive before this date)
# so it is intentionally skipped — it gets no BX cell, matching the source.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New date header (2025-11-07) --------------------------------------
$ws.Cells.Item(1, 76).Value2 = 45968

# --- Per-player attendance for the new date ----------------------------
$ws.Cells.Item(2, 76).Value = "P"
$ws.Cells.Item(3, 76).Value = "R"
$ws.Cells.Item(4, 76).Value = "P"
$ws.Cells.Item(5, 76).Value = "B"
$ws.Cells.Item(6, 76).Value = "B"
$ws.Cells.Item(7, 76).Value = "P"
$ws.Cells.Item(8, 76).Value = "RH"
$ws.Cells.Item(9, 76).Value = "P"
$ws.Cells.Item(10, 76).Value = "P"
$ws.Cells.Item(11, 76).Value = "P"
# row 12 intentionally left blank (no BX cell in source)
$ws.Cells.Item(13, 76).Value = "B"
$ws.Cells.Item(14, 76).Value = "M"
$ws.Cells.Item(15, 76).Value = "P"
$ws.Cells.Item(16, 76).Value = "P"
$ws.Cells.Item(17, 76).Value = "P"
$ws.Cells.Item(18, 76).Value = "B"
$ws.Cells.Item(19, 76).Value = "P"
$ws.Cells.Item(20, 76).Value = "P"
# row 21 stays empty (no mark recorded for this date)
$ws.Cells.Item(22, 76).Value = "P"
$ws.Cells.Item(23, 76).Value = "RH"
$ws.Cells.Item(24, 76).Value = "P"
$ws.Cells.Item(25, 76).Value = "P"
$ws.Cells.Item(26, 76).Value = "P"
$ws.Cells.Item(27, 76).Value = "P"
$ws.Cells.Item(28, 76).Value = "P"
$ws.Cells.Item(29, 76).Value = "B"

# --- Small correction noticed on row 3 (BV3: R -> P) --------------------
$ws.Cells.Item(3, 74).Value = "P"

# --- Copy formatting from BW into the new BX column ---------------------
# Done in two pieces so row 12 (which has no BW entry) is left untouched.
# This runs AFTER the values are set so the earlier formula recalculation
# isn't short-circuited by the paste operation.
$ws.Range("BW1:BW11").Copy()
$ws.Range("BX1:BX11").PasteSpecial(-4122)

$ws.Range("BW13:BW29").Copy()
$ws.Range("BX13:BX29").PasteSpecial(-4122)

# --- Refresh the view state: new frozen-pane scroll + selected cell -----
$ws.Range("BX1").Select()
$ws.Range("CA19").Select()
